# Update countries & provincias Spain
#
# The live "paises" COVID dashboard was re-pulled (11:11 -> 12:28 timestamp).
# Three countries jumped rank in the Casos-totales sort, which moves their
# name into a neighbouring row, and the per-row totals for every row whose
# data changed (including the rows that shifted) are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country-name (column A) updates caused by the re-sort ---
$ws.Range("A1").Value = 'Datos actualizados a 24 de Octubre de 2020 a las 12:28'
$ws.Range("A41").Value = 'Emiratos Arabes Unidos'
$ws.Range("A42").Value = 'Republica Dominicana'
$ws.Range("A141").Value = 'Letonia'
$ws.Range("A142").Value = 'Aruba'
$ws.Range("A143").Value = 'Estonia'
$ws.Range("A144").Value = 'Islandia'
$ws.Range("A145").Value = 'Mayotte'
$ws.Range("A216").Value = 'Islas Malvinas'
$ws.Range("A217").Value = 'Montserrat'

# --- Refreshed statistics (columns B:H) for the affected rows ---
# row 21
$ws.Range("B21").Value = 397507
$ws.Range("C21").Value = 1094
$ws.Range("D21").Value = 313563
$ws.Range("E21").Value = 78164
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 19
$ws.Range("H21").Value = 5780
# row 34
$ws.Range("B34").Value = 205793
$ws.Range("C34").Value = 4761
$ws.Range("D34").Value = 147932
$ws.Range("E34").Value = 51543
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 73
$ws.Range("H34").Value = 6318
# row 41
$ws.Range("B41").Value = 123764
$ws.Range("C41").Value = 1491
$ws.Range("D41").Value = 116894
$ws.Range("E41").Value = 6395
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 475
# row 42
$ws.Range("B42").Value = 122873
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 100920
$ws.Range("E42").Value = 19741
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 2212
# row 50
$ws.Range("B50").Value = 103653
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 55800
$ws.Range("E50").Value = 45782
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 2071
# row 71
$ws.Range("B71").Value = 54374
$ws.Range("C71").Value = 990
$ws.Range("D71").Value = 29965
$ws.Range("E71").Value = 23619
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 16
$ws.Range("H71").Value = 790
# row 92
$ws.Range("B92").Value = 25742
$ws.Range("C92").Value = 1228
$ws.Range("D92").Value = 16555
$ws.Range("E92").Value = 8966
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 7
$ws.Range("H92").Value = 221
# row 133
$ws.Range("B133").Value = 5290
$ws.Range("C133").Value = 5
$ws.Range("D133").Value = 5041
$ws.Range("E133").Value = 144
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 105
# row 138
$ws.Range("B138").Value = 5079
$ws.Range("C138").Value = 5
$ws.Range("D138").Value = 4962
$ws.Range("E138").Value = 34
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 83
# row 141
$ws.Range("B141").Value = 4467
$ws.Range("C141").Value = 259
$ws.Range("D141").Value = 1357
$ws.Range("E141").Value = 3056
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 4
$ws.Range("H141").Value = 54
# row 142
$ws.Range("B142").Value = 4401
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 4160
$ws.Range("E142").Value = 205
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 36
# row 143
$ws.Range("B143").Value = 4351
$ws.Range("C143").Value = 52
$ws.Range("D143").Value = 3441
$ws.Range("E143").Value = 837
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 73
# row 144
$ws.Range("B144").Value = 4308
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 3187
$ws.Range("E144").Value = 1110
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 11
# row 145
$ws.Range("B145").Value = 4276
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 2964
$ws.Range("E145").Value = 1268
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 44
# row 152
$ws.Range("B152").Value = 3314
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 1882
$ws.Range("E152").Value = 1407
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 25
# row 216
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
# row 217
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 1

Write-Output "Applied country reorder + stats refresh"
